# Apply "Upload new version with timestamp" edit:
# Fill in the first data row (row 7) of the "نواقص الأصناف" (items shortage)
# report with a concrete item: COLOVATIL 30 F.C. TABS, along with its
# balance, order limit, price, selling price and transaction count.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Merged range C7:G7 -> item name
$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value = "COLOVATIL 30 F.C. TABS"

# Merged range H7:K7 -> current balance
$ws.Range("H7").NumberFormat = "@"
$ws.Range("H7").Value = "0:1"

# Merged range L7:M7 -> order limit
$ws.Range("L7").Value = "1"

# Merged range N7:O7 -> price
$ws.Range("N7").NumberFormat = "@"
$ws.Range("N7").Value = "63.00"

# P7 -> selling price
$ws.Range("P7").Value = "63.0000"

# Q7 -> number of transactions
$ws.Range("Q7").NumberFormat = "@"
$ws.Range("Q7").Value = "1:0"
